# Update countries & provincias Spain
# Applies the data refresh captured in the commit: swap the Vietnam /
# Montenegro rows (they traded places in the source ranking), refresh a
# handful of per-country case/death counters, and bump the "last updated"
# timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the Vietnam / Montenegro country labels (rows 142 & 143) ---
# Row 142 used to be Montenegro, row 143 used to be Vietnam; after the
# refresh they swap places and each picks up new stats.
$ws.Range("A142").Value = "Vietnam"
$ws.Range("A143").Value = "Montenegro"

# --- Row 52: Noruega ---
$ws.Range("E52").Value = 7984
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 233

# --- Row 61: Finlandia ---
$ws.Range("E61").Value = 1080
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 300

# --- Row 93: Eslovenia ---
$ws.Range("D93").Value = 274
$ws.Range("E93").Value = 1088

# --- Row 124: Malta ---
$ws.Range("B124").Value = 558
$ws.Range("C124").Value = 5
$ws.Range("D124").Value = 456
$ws.Range("E124").Value = 96

# --- Row 142: now Vietnam ---
$ws.Range("C142").Value = 4
$ws.Range("D142").Value = 263
$ws.Range("E142").Value = 61
$ws.Range("H142").Value = 0

# --- Row 143: now Montenegro ---
$ws.Range("B143").Value = 324
$ws.Range("D143").Value = 311
$ws.Range("E143").Value = 4
$ws.Range("H143").Value = 9

# --- Update "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 13:35"
